# Apply scheduled-runner profit/price updates across the Sagittarius_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 142
$ws.Range("I33").Value = 141
$ws.Range("K33").Value = 141
$ws.Range("M33").Value = 88
$ws.Range("H86").Value = 1708.5
$ws.Range("I86").Value = 1383.1666
$ws.Range("K86").Value = 1383.1666
$ws.Range("M86").Value = -260.1666
$ws.Range("H89").Value = 1708.5
$ws.Range("I89").Value = 1383.1666
$ws.Range("K89").Value = 6915.833000000001
$ws.Range("M89").Value = -1299.833000000001
$ws.Range("H92").Value = 229.33333
$ws.Range("I92").Value = 229.33333
$ws.Range("K92").Value = 229.33333
$ws.Range("M92").Value = 1018.66667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1351.5625
$ws.Range("I74").Value = 1331.6364
$ws.Range("J74").Value = 1395.4
$ws.Range("K74").Value = 1331.6364
$ws.Range("L74").Value = 1395.4
$ws.Range("M74").Value = -457.6364000000001
$ws.Range("N74").Value = -3143.4
$ws.Range("H77").Value = 1351.5625
$ws.Range("I77").Value = 1331.6364
$ws.Range("J77").Value = 1395.4
$ws.Range("K77").Value = 6658.182000000001
$ws.Range("L77").Value = 6977
$ws.Range("M77").Value = -2290.182000000001
$ws.Range("N77").Value = -15713
$ws.Range("H122").Value = 2358.25
$ws.Range("I122").Value = 1644.3334
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 4933.0002
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -2483.0002
$ws.Range("N122").Value = -18400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9999.5
$ws.Range("I20").Value = 9999
$ws.Range("K20").Value = 9999
$ws.Range("M20").Value = -9752
$ws.Range("H107").Value = 1423.3
$ws.Range("I107").Value = 1373.3334
$ws.Range("K107").Value = 1373.3334
$ws.Range("M107").Value = 546.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3857.7144
$ws.Range("J16").Value = 2750
$ws.Range("L16").Value = 2750
$ws.Range("N16").Value = -3324
$ws.Range("H31").Value = 2166.25
$ws.Range("I31").Value = 2166.25
$ws.Range("K31").Value = 2166.25
$ws.Range("M31").Value = -1871.25
$ws.Range("H34").Value = 2166.25
$ws.Range("I34").Value = 2166.25
$ws.Range("K34").Value = 2166.25
$ws.Range("M34").Value = -1964.25
$ws.Range("H107").Value = 1286.1818
$ws.Range("I107").Value = 912.1429000000001
$ws.Range("K107").Value = 912.1429000000001
$ws.Range("M107").Value = 1007.8571
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").ClearContents()
$ws.Range("N108").Value = 0
$ws.Range("H113").Value = 3857.7144
$ws.Range("J113").Value = 2750
$ws.Range("L113").Value = 2750
$ws.Range("N113").Value = -7090
$ws.Range("H122").Value = 950
$ws.Range("I122").Value = 950
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2850
$ws.Range("L122").ClearContents()
$ws.Range("N122").Value = 0
$ws.Range("M122").Value = -400
$ws.Range("H132").Value = 4332.3335
$ws.Range("I132").Value = 4500
$ws.Range("K132").Value = 13500
$ws.Range("M132").Value = -10970

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1576.5
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 1576.5
$ws.Range("K34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("M34").Value = 4729.5
$ws.Range("N34").Value = -4897.5
$ws.Range("H113").Value = 1649.5
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 44985
$ws.Range("J15").Value = 44985
$ws.Range("L15").Value = 44985
$ws.Range("N15").Value = -45561
$ws.Range("H70").Value = 4975
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 4975
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H81").Value = 44985
$ws.Range("J81").Value = 44985
$ws.Range("L81").Value = 44985
$ws.Range("N81").Value = -46981
$ws.Range("H84").Value = 44985
$ws.Range("J84").Value = 44985
$ws.Range("L84").Value = 134955
$ws.Range("N84").Value = -144939
$ws.Range("H132").Value = 6359.5356
$ws.Range("I132").Value = 6522.92
$ws.Range("K132").Value = 19568.76
$ws.Range("M132").Value = -17038.76
$ws.Range("H134").Value = 23658.8
$ws.Range("J134").Value = 23658.8
$ws.Range("L134").Value = 70976.39999999999
$ws.Range("N134").Value = -76046.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4166.6665
$ws.Range("I22").Value = 4600
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 4600
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -4305
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 4166.6665
$ws.Range("I27").Value = 4600
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 4600
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -4493
$ws.Range("N27").Value = -2214
$ws.Range("H46").Value = 1481.6923
$ws.Range("I46").Value = 1373.375
$ws.Range("J46").Value = 1529.8334
$ws.Range("K46").Value = 1373.375
$ws.Range("L46").Value = 1529.8334
$ws.Range("M46").Value = -1185.375
$ws.Range("N46").Value = -1905.8334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 49595
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 49595
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H126").Value = 3050.4285
$ws.Range("I126").Value = 1984.6666
$ws.Range("K126").Value = 5953.9998
$ws.Range("M126").Value = -3483.9998
$ws.Range("H136").Value = 2021.2
$ws.Range("I136").Value = 1369
$ws.Range("J136").Value = 2999.5
$ws.Range("K136").Value = 4107
$ws.Range("L136").Value = 8998.5
$ws.Range("M136").Value = -1557
$ws.Range("N136").Value = -14098.5
